$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value that was updated for every
# data row (rows 2-171) from 2023-09-03 (45172) to 2023-09-06 (45175).
$ws.Range("C2:C171").Value = 45175
